# Motion planning writeup.docx - apply the commit's edits:
#   1. Remove the stray "_GoBack" bookmark from the very start of the document
#      (it was left over at the top of paragraph 1).
#   2. Append two new paragraphs (separated by blank paragraphs) at the end of
#      the body, right before the trailing blank paragraph / sectPr, and leave
#      a fresh "_GoBack" bookmark collapsed at the very end of the new text
#      (mirroring where Word drops it after the last edit position).

$d = $word.ActiveDocument

# --- 1. Remove the old "_GoBack" bookmark -----------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- 2. Insert the new paragraphs just before the final (empty) paragraph ---
$lastParaIndex = $d.Paragraphs.Count
$insertionPoint = $d.Paragraphs.Item($lastParaIndex).Range.Duplicate
$insertionPoint.Collapse(1)   # wdCollapseStart - insert before the trailing blank paragraph

$newContentXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:p><w:r><w:t xml:space="preserve">I modified the location of the goal by giving a different default and this way I could show it would work with other locations and would also adapt to give goal locations. </w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">I used many resources to complete this project including the class notes, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>github</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> examples, student forum and even the slack channel.  There are pieces of code modified based on these resources. </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$insertionPoint.InsertXML($newContentXml)

Write-Output "Paragraph count after edit: $($d.Paragraphs.Count)"
